# Add a new "MAY-22" worksheet after "APR-22" and populate it with the
# May daily-track entries, mirroring the structure of the other monthly
# sheets (JAN-22 / FEB-22 / MAR-22 / APR-22).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reference cells (already-saved styles) that we reuse via copy/paste of
# formats only, so that the new cells land on the very same style (xf)
# index as is already used elsewhere in the workbook instead of growing
# the style table with duplicate entries.
# ---------------------------------------------------------------------
$dateStyleSrc    = $wb.Worksheets.Item("JAN-22").Range("B30")   # date, no border
$wrapStyleSrc    = $wb.Worksheets.Item("APR-22").Range("D9")    # wrap text, no border/fill
$percentStyleSrc = $wb.Worksheets.Item("APR-22").Range("E9")    # percentage, no border/fill

# ---------------------------------------------------------------------
# Create the new sheet right after the current last sheet (APR-22) and
# rename it, then make it the active / selected tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "MAY-22"

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value2 = "No"
$ws.Range("B1").Value2 = "Date"
$ws.Range("C1").Value2 = "Application"
$ws.Range("D1").Value2 = "Task"
$ws.Range("E1").Value2 = "% of completion"
$ws.Range("F1").Value2 = "Status"
$ws.Range("G1").Value2 = "Comments"

# ---------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = 44683
$ws.Range("C2").Value2 = "RPA GSS"
$ws.Range("D2").Value2 = "1. Service Order Parts Number task has been executed to extract the Parts number and uploading the csv files for SSC1 and SSC3, whereas the other SSCs are work in progress"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = "Completed"

$ws.Range("D3").Value2 = "2. Working on Exe conversion for the GL, P&L, and Management Template works etc and it is work in progress"
$ws.Range("E3").Value2 = 0.6
$ws.Range("F3").Value2 = "WIP"

# ---------------------------------------------------------------------
# Formatting - reuse existing styles where possible via format-only paste
# ---------------------------------------------------------------------
$dateStyleSrc.Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wrapStyleSrc.Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null

$percentStyleSrc.Copy() | Out-Null
$ws.Range("E2:E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# A1 uses a (new) number format of "mmm-yy" even though it holds text.
$ws.Range("A1").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------
# Column widths (best-fit-like, matching the other monthly sheets)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334
$ws.Columns.Item(4).ColumnWidth = 115.5
$ws.Columns.Item(7).ColumnWidth = 9.0

# ---------------------------------------------------------------------
# Row height for the wrapped second row
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 28.8

# ---------------------------------------------------------------------
# View / selection state: MAY-22 becomes the active tab, with C1 scrolled
# into view and D2:D3 selected; APR-22 loses its tabSelected flag and its
# own selection moves from D33 to D32.
# ---------------------------------------------------------------------
$apr = $wb.Worksheets.Item("APR-22")
$apr.Activate()
$apr.Range("D32").Select() | Out-Null

$ws.Activate()
$ws.Range("C1").Select() | Out-Null
$ws.Range("D2:D3").Select() | Out-Null
